$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-20 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-09-21 Sunday", 2) | Out-Null
$d.Content.Find.Execute("35×88=3080", $true, $true, $false, $false, $false, $true, 1, $false, "51×86=4386", 2) | Out-Null
$d.Content.Find.Execute("32×62=1984", $true, $true, $false, $false, $false, $true, 1, $false, "71×41=2911", 2) | Out-Null
$d.Content.Find.Execute("18×79=1422", $true, $true, $false, $false, $false, $true, 1, $false, "31×66=2046", 2) | Out-Null
$d.Content.Find.Execute("96×56=5376", $true, $true, $false, $false, $false, $true, 1, $false, "86×17=1462", 2) | Out-Null
$d.Content.Find.Execute("45×91=4095", $true, $true, $false, $false, $false, $true, 1, $false, "43×23=989", 2) | Out-Null
$d.Content.Find.Execute("91×90=8190", $true, $true, $false, $false, $false, $true, 1, $false, "82×49=4018", 2) | Out-Null
$d.Content.Find.Execute("92×97=8924", $true, $true, $false, $false, $false, $true, 1, $false, "50×13=650", 2) | Out-Null
$d.Content.Find.Execute("72×37=2664", $true, $true, $false, $false, $false, $true, 1, $false, "93×12=1116", 2) | Out-Null
$d.Content.Find.Execute("57×69=3933", $true, $true, $false, $false, $false, $true, 1, $false, "83×64=5312", 2) | Out-Null
$d.Content.Find.Execute("57×17=969", $true, $true, $false, $false, $false, $true, 1, $false, "91×94=8554", 2) | Out-Null
$d.Content.Find.Execute("76×20=1520", $true, $true, $false, $false, $false, $true, 1, $false, "45×38=1710", 2) | Out-Null
$d.Content.Find.Execute("42×70=2940", $true, $true, $false, $false, $false, $true, 1, $false, "24×91=2184", 2) | Out-Null
$d.Content.Find.Execute("35×91=3185", $true, $true, $false, $false, $false, $true, 1, $false, "41×59=2419", 2) | Out-Null
$d.Content.Find.Execute("20×65=1300", $true, $true, $false, $false, $false, $true, 1, $false, "53×65=3445", 2) | Out-Null
$d.Content.Find.Execute("43×78=3354", $true, $true, $false, $false, $false, $true, 1, $false, "66×20=1320", 2) | Out-Null
$d.Content.Find.Execute("61×56=3416", $true, $true, $false, $false, $false, $true, 1, $false, "87×55=4785", 2) | Out-Null
$d.Content.Find.Execute("44×23=1012", $true, $true, $false, $false, $false, $true, 1, $false, "15×34=510", 2) | Out-Null
$d.Content.Find.Execute("18×23=414", $true, $true, $false, $false, $false, $true, 1, $false, "69×16=1104", 2) | Out-Null
$d.Content.Find.Execute("46×88=4048", $true, $true, $false, $false, $false, $true, 1, $false, "39×19=741", 2) | Out-Null
$d.Content.Find.Execute("86×99=8514", $true, $true, $false, $false, $false, $true, 1, $false, "88×89=7832", 2) | Out-Null
$d.Content.Find.Execute("66×68=4488", $true, $true, $false, $false, $false, $true, 1, $false, "81×58=4698", 2) | Out-Null
$d.Content.Find.Execute("18×42=756", $true, $true, $false, $false, $false, $true, 1, $false, "31×35=1085", 2) | Out-Null
$d.Content.Find.Execute("47×42=1974", $true, $true, $false, $false, $false, $true, 1, $false, "84×48=4032", 2) | Out-Null
$d.Content.Find.Execute("37×87=3219", $true, $true, $false, $false, $false, $true, 1, $false, "52×34=1768", 2) | Out-Null
$d.Content.Find.Execute("42×96=4032", $true, $true, $false, $false, $false, $true, 1, $false, "73×62=4526", 2) | Out-Null
